$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Add a bottom paragraph border (horizontal separator) under the third
#    paragraph ("It also contains a header and a couple of paragraphs. ...
#    hyperlink.").
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
$p3.Borders.Item(-3).LineStyle = 1      # wdLineStyleSingle
$p3.Borders.Item(-3).LineWidth = 3      # -> w:sz="6" (0.75pt)
$p3.Borders.Item(-3).ColorIndex = 0     # wdAuto -> w:color="auto"
$p3.Borders.DistanceFromBottom = 1      # -> w:space="1"

# ---------------------------------------------------------------------------
# 2) Append two new paragraphs after paragraph 3, with mixed run formatting,
#    using WordOpenXML package fragments via Range.InsertXML so that the
#    exact rPr (b/bCs, i/iCs, u) combinations are reproduced faithfully.
# ---------------------------------------------------------------------------
$p3.Range.InsertParagraphAfter()

$pkgHeader = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$pkgFooter = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

# --- Paragraph A: "The test also tests partial formats. It also tests formats across"
#     (paragraph mark itself is bold)
$paraA = '<w:body><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>The test also tests pa</w:t></w:r>' +
  '<w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>rtia</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>l f</w:t></w:r>' +
  '<w:r><w:rPr><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr><w:t>o</w:t></w:r>' +
  '<w:r><w:rPr><w:i/><w:iCs/><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr><w:t>r</w:t></w:r>' +
  '<w:r><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr><w:t>m</w:t></w:r>' +
  '<w:r><w:rPr><w:i/><w:iCs/><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr><w:t>a</w:t></w:r>' +
  '<w:r><w:rPr><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr><w:t>t</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">s. It also tests formats </w:t></w:r>' +
  '<w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>across</w:t></w:r>' +
  '</w:p></w:body>'

$p4 = $d.Paragraphs.Item(4)
$r4 = $p4.Range
$r4.InsertXML($pkgHeader + $paraA + $pkgFooter)

# --- Paragraph B: "paragraph boundaries."
$paraB = '<w:body><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>paragraph</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> boundaries.</w:t></w:r>' +
  '</w:p></w:body>'

$p5 = $d.Paragraphs.Item(5)
$r5 = $p5.Range
$r5.InsertXML($pkgHeader + $paraB + $pkgFooter)

# InsertXML leaves a trailing empty paragraph behind (the split remainder of
# the paragraph mark it was inserted in front of). Merge it away so the
# document ends cleanly with paragraph B, right before the sectPr.
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$prevPara = $d.Paragraphs.Item($count - 1)
if ($lastPara.Range.Start -eq $lastPara.Range.End - 1) {
  $mergeRange = $d.Range($prevPara.Range.End - 1, $lastPara.Range.End)
  $mergeRange.Delete()
}

Write-Host "Final paragraph count:" $d.Paragraphs.Count
